$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 31, shifting existing rows 31-155 down to 33-157
$ws.Range("A31:A32").EntireRow.Insert()

# Populate new row 31
$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C31").Value = "Ñuble"
$ws.Range("D31").Value = 44676
$ws.Range("D31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E31").Value = 16
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100101
$ws.Range("H31").Value = "Berries"
$ws.Range("I31").Value = 100101007
$ws.Range("J31").Value = "Kiwi"
$ws.Range("K31").Value = "Hayward"
$ws.Range("L31").Value = "Especial"
$ws.Range("M31").Value = 30
$ws.Range("N31").Value = 11000
$ws.Range("O31").Value = 11000
$ws.Range("P31").Value = 11000
$ws.Range("Q31").Value = "`$/bandeja 18 kilos"
$ws.Range("R31").Value = "Provincia de Curicó"
$ws.Range("S31").Value = 611
$ws.Range("T31").Value = 18

# Populate new row 32
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44676
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100101
$ws.Range("H32").Value = "Berries"
$ws.Range("I32").Value = 100101007
$ws.Range("J32").Value = "Kiwi"
$ws.Range("K32").Value = "Hayward"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 60
$ws.Range("N32").Value = 9000
$ws.Range("O32").Value = 10000
$ws.Range("P32").Value = 9500
$ws.Range("Q32").Value = "`$/bandeja 18 kilos"
$ws.Range("R32").Value = "Provincia de Curicó"
$ws.Range("S32").Value = 528
$ws.Range("T32").Value = 18
